$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need to be forced to
# Text format first, otherwise Excel auto-converts them to floating point
# numbers (losing exact formatting like trailing zeros, e.g. "2.40" -> 2.4,
# or introducing precision noise, e.g. "112.68" -> 112.68000000000001).
$textRefs = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D16", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D45", "D47", "D49", "D50", "D51")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.647.10'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '2.228.64'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '112.68'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').Value = '293.55'
$ws.Range('E6').Value = '  +10.39%  '
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').Value = '0.602'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').Value = '43.71'
$ws.Range('E10').Value = '  -5.61%  '
$ws.Range('D11').Value = '0.0916'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').Value = '54.52'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = '8.66'
$ws.Range('E13').Value = '  -6.08%  '
$ws.Range('E14').Value = '  +19.65%  '
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').Value = '14.99'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').Value = '2.562.53'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '2.226.08'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').Value = '42.498.56'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').Value = '7.22'
$ws.Range('E20').Value = '  +7.34%  '
$ws.Range('D21').Value = '0.0000106'
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('D22').Value = '73.52'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('E23').Value = '  +15.19%  '
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').Value = '236.07'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').Value = '8.95'
$ws.Range('E26').Value = '  -3.85%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').Value = '11.46'
$ws.Range('E28').Value = '  -8.69%  '
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').Value = '174.92'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').Value = '37.47'
$ws.Range('E31').Value = '  -7.88%  '
$ws.Range('D32').Value = '3.13'
$ws.Range('E32').Value = '  -5.07%  '
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('D34').Value = '0.0885'
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').Value = '5.67'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('D36').Value = '5.04'
$ws.Range('E36').Value = '  +9.72%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '0.126'
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('D41').Value = '2.40'
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('D42').Value = '71.67'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '12.35'
$ws.Range('E45').Value = '  -7.75%  '
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').Value = '5.38'
$ws.Range('E47').Value = '  -4.69%  '
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').Value = '8.44'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.64'
$ws.Range('E50').Value = '  +5.01%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '101.33'
$ws.Range('E51').Value = '  +1.20%  '
